# Daily automated data refresh for the EPEX Spot workbook.
#  - "Prix Spot": append a new day column (CP) with the 15-sep hourly prices.
#  - "Gaz" / "CO2": append the two most recent daily quotes (2025-09-13 / 2025-09-14).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": new column CP ("15-sep"), one header cell + 24 hourly rows.
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Header cell - copy the previous day's header (for formatting/style: bold,
# bordered, centered) then overwrite its text with the new date label.
$wsPrix.Range("CO1").Copy($wsPrix.Range("CP1"))
$wsPrix.Range("CP1").Value = "15-sep"

# Hourly values for the new day, row 2 (00 - 01) through row 25 (23 - 24).
$cpValues = [ordered]@{
    2  = 0
    3  = -0.01
    4  = -0.01
    5  = -0.01
    6  = -0.01
    7  = -0.01
    8  = 0
    9  = 2.18
    10 = 19.69
    11 = 0.9
    12 = 0
    13 = -0.01
    14 = -0.04
    15 = -0.02
    16 = -0.06
    17 = -0.01
    18 = 0
    19 = 0.43
    20 = 6.03
    21 = 28.19
    22 = 9.369999999999999
    23 = 5.16
    24 = 5.17
    25 = 4.29
}

foreach ($row in $cpValues.Keys) {
    # Copy formatting from the same row's previous-day cell, then set the value.
    $wsPrix.Range("CO$row").Copy($wsPrix.Range("CP$row"))
    $wsPrix.Range("CP$row").Value = $cpValues[$row]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append rows 91 (2025-09-13) and 92 (2025-09-14), value 32.2.
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

$gazRows = [ordered]@{
    91 = "2025-09-13"
    92 = "2025-09-14"
}

foreach ($row in $gazRows.Keys) {
    # Force the date column to Text first so the "yyyy-mm-dd" string is kept
    # verbatim instead of being auto-parsed into a date serial number, then
    # drop the explicit number-format again so the cell stays unstyled like
    # its neighbours.
    $wsGaz.Range("A$row").NumberFormat = "@"
    $wsGaz.Range("A$row").Value = $gazRows[$row]
    $wsGaz.Range("A$row").ClearFormats()

    $wsGaz.Range("B$row").Value = 32.2
}

# ---------------------------------------------------------------------------
# Sheet "CO2": append rows 91 (2025-09-13) and 92 (2025-09-14), value 75.47.
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$co2Rows = [ordered]@{
    91 = "2025-09-13"
    92 = "2025-09-14"
}

foreach ($row in $co2Rows.Keys) {
    $wsCo2.Range("A$row").NumberFormat = "@"
    $wsCo2.Range("A$row").Value = $co2Rows[$row]
    $wsCo2.Range("A$row").ClearFormats()

    $wsCo2.Range("B$row").Value = 75.47
}
